$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 12 (cohort_year 2021, period_index 4): num_customers 62 -> 63, retention_rate recalculated
$ws.Range("C12").Value = 63
$ws.Range("E12").Value = 0.1503579952267303

# Row 16 (cohort_year 2022, period_index 3): num_customers 51 -> 53, retention_rate recalculated
$ws.Range("C16").Value = 53
$ws.Range("E16").Value = 0.2746113989637305

# Row 21 (cohort_year 2024, period_index 1): num_customers 91 -> 96, retention_rate recalculated
$ws.Range("C21").Value = 96
$ws.Range("E21").Value = 0.4660194174757282

# Row 22 (cohort_year 2025, period_index 0): num_customers 16 -> 20, cohort_size 16 -> 20, retention_rate stays 1
$ws.Range("C22").Value = 20
$ws.Range("D22").Value = 20
$ws.Range("E22").Value = 1
